$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new "Markertypes" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "S-Matrix"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Markertypes"

# --- Markertype[q] table (columns A:B) ---
$ws2.Range("A1").Value = "Markertype[q]"
$ws2.Range("B1").Value = "Entry"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Markername"

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "Upos"

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "Vpos"

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "SizeU"

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "SizeV"

$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "StepU"

$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = "StepV"

$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "PointsU"

$ws2.Range("A10").Value = 8
$ws2.Range("B10").Value = "PointsV"

$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = "MarkOffsetU"

$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = "MarkOffsetV"

$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = "MarkplaceU"

$ws2.Range("A14").Value = 12
$ws2.Range("B14").Value = "MarkplaceV"

$ws2.Range("A15").Value = 13
$ws2.Range("B15").Value = "Profile min"

$ws2.Range("A16").Value = 14
$ws2.Range("B16").Value = "Profile max"

$ws2.Range("A17").Value = 15
$ws2.Range("B17").Value = "ContrastLow"

$ws2.Range("A18").Value = 16
$ws2.Range("B18").Value = "ContrastHigh"

$ws2.Range("A19").Value = 17
$ws2.Range("B19").Value = "Threshold"

$ws2.Range("A20").Value = 18
$ws2.Range("B20").Value = "Entry positionlist"

# --- WFAlignprocedures[q] table (columns D:E) ---
$ws2.Range("D1").Value = "WFAlignprocedures[q]"

$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = "Procedurename"

$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = "entries (inc. log)"

$ws2.Range("D4").Value = 2
$ws2.Range("E4").Value = "1st entry"

$ws2.Range("D5").Value = 3
$ws2.Range("E5").Value = "2nd entry"

# D6 was typed with a leading apostrophe in the source workbook, which is why
# it carries the quotePrefix style while the identical D7 text does not.
$ws2.Range("D6").Value = "'.."
$ws2.Range("E6").Value = "etc"

$ws2.Range("D7").Value = ".."

$ws2.Range("D8").Value = "last"
$ws2.Range("E8").Value = "log"

# --- Cosmetics: best-fit the two label columns like Excel would ---
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null

# --- Selection / active sheet state ---
$ws1.Range("C54").Select() | Out-Null
$ws2.Range("E4").Select() | Out-Null

$ws2.Activate() | Out-Null
